# Automatische test-sync: 2025-06-26 23:51:50
#
# Appends the 15th test-mail log entry (row 47) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover the new row, and
# updates the "Dashboard" summary table so the "Overig" category count
# (now 3) outranks "Offerte / Prijsaanvraag" (2), swapping their rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Logs sheet: append row 47
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A47").Value = "Ik wil dat je dit regelt met support."
$logs.Range("B47").Value = "mailmind.test@zohomail.eu"
$logs.Range("C47").Value = "Testmail #15: Ik wil dat je dit regelt met support."
$logs.Range("D47").Value = "Overig"

$antwoord47 = @"
Beste klant,
Bedankt voor uw e-mail. Om uw verzoek efficiënt te kunnen afhandelen, zou ik graag wat meer informatie van u willen ontvangen. Kunt u mij alstublieft meer details geven over het specifieke probleem waarmee u hulp nodig heeft? Op die manier kunnen wij u beter van dienst zijn.
Met vriendelijke groet,
[Naam]  
E-mailassistent  
[Bedrijfsnaam]
"@
$logs.Range("E47").Value = $antwoord47

$logs.Range("F47").Value = "2025-06-26 23:51:43"
$logs.Range("G47").Value = "Ja"
$logs.Range("H47").Value = "Nee"
$logs.Range("I47").Value = "Ja"

# Re-fit the row height (writing the multi-line E47 text auto-expands
# it); AutoFit restores the sheet's default row metrics like the rest
# of the log rows, which were never given an explicit height.
$logs.Rows.Item(47).AutoFit()

# ---------------------------------------------------------------
# 2) Extend conditional formatting sqref from row 46 to row 47
#    (one ModifyAppliesToRange per block keeps all cfRules in that
#    block in sync, same as Excel's UI "Applies to" edit).
# ---------------------------------------------------------------
$logs.Range("D2:D46").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D47"))
$logs.Range("G2:G46").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G47"))
$logs.Range("H2:H46").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H47"))
$logs.Range("I2:I46").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I47"))

# ---------------------------------------------------------------
# 3) Dashboard sheet: "Overig" now counts 3 (was 2), overtaking
#    "Offerte / Prijsaanvraag" (2) in the sorted summary table, so
#    rows 6 and 7 swap places.
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Overig"
$dash.Range("B6").Value = 3
$dash.Range("A7").Value = "Offerte / Prijsaanvraag"
$dash.Range("B7").Value = 2

Write-Host "Row 47 added to Logs, conditional formatting extended, Dashboard rows 6-7 swapped."
